$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date value for rows 2 through 15
# from 45170 (2023-09-01) to 45174 (2023-09-05)
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
